# Legacy GSC export data refresh: the underlying export moved forward by one
# day, so the oldest data row (2025-09-12, which had no video-index data yet)
# drops off the top of the "Chart" table and every subsequent row shifts up
# by one. Deleting the first data row reproduces that shift natively (dates,
# values, and the trailing row all cascade up correctly) and shrinks the
# sheet's used range from A1:D89 to A1:D88, matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
